# Generate Report for Handoff
#
# The localization-status report moves from "In Translation" to
# "Ready for handoff": the Overview sheet's per-language status columns and
# each per-language sheet's Status column get the new text, and the
# "Latest Handoff Datetime" / "Latest HO Xliff Generate Date" timestamps are
# refreshed to the moment the handoff report was produced. Excel re-flows the
# (now wider) Status columns to fit the longer text.

$wb = $excel.ActiveWorkbook

$oldStatus = "In Translation"
$newStatus = "Ready for handoff"

# --- Overview sheet: zh-cn / de-de status columns (E, F) + HO xliff date (G) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("G2").Value = "2016-08-25 02:38:49"

# --- zh-cn sheet: Status (C) + Latest Handoff Datetime (H) ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("H2").Value = "2016-08-25 02:38:45"

# --- de-de sheet: Status (C) + Latest Handoff Datetime (H) ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus
$dede.Range("H2").Value = "2016-08-25 02:38:49"

# --- Re-flow the Status columns now that "Ready for handoff" is wider than
#     "In Translation" ---
$overview.Columns.Item(5).ColumnWidth = 16.333333333333336
$overview.Columns.Item(6).ColumnWidth = 16.333333333333336
$zhcn.Columns.Item(3).ColumnWidth = 16.333333333333336
$dede.Columns.Item(3).ColumnWidth = 16.333333333333336

Write-Host "Updated status to '$newStatus' and refreshed handoff timestamps"
